$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.409.55'
$ws.Range("E2").Value = '  -1.18%  '
$ws.Range("D3").Value = '3.065.42'
$ws.Range("E3").Value = '  -2.52%  '
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '588.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.76'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.12%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("E8").Value = '  +0.74%  '
$ws.Range("D9").Value = '3.062.00'
$ws.Range("E9").Value = '  -2.31%  '
$ws.Range("E10").Value = '  -4.13%  '
$ws.Range("E11").Value = '  -1.53%  '
$ws.Range("E12").Value = '  -1.57%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '36.96'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.49%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000237'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.19%  '
$ws.Range("E15").Value = '  -1.89%  '
$ws.Range("D16").Value = '3.574.11'
$ws.Range("E16").Value = '  -2.48%  '
$ws.Range("D17").Value = '63.433.42'
$ws.Range("E17").Value = '  -0.81%  '
$ws.Range("E18").Value = '  -2.00%  '
$ws.Range("D19").Value = '3.061.30'
$ws.Range("E19").Value = '  -2.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '471.78'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.87%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.32'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.51'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.95%  '
$ws.Range("E24").Value = '  +1.21%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '80.58'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.55%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.79'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.56%  '
$ws.Range("E27").Value = '  +4.53%  '
$ws.Range("E28").Value = '  -0.25%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.48'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.83%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.31%  '
$ws.Range("E31").Value = '  -2.06%  '
$ws.Range("E32").Value = '  -3.05%  '
$ws.Range("E33").Value = '  -1.77%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '27.16'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.90%  '
$ws.Range("D35").Value = '0.0₃0818'
$ws.Range("E35").Value = '  -4.45%  '
$ws.Range("E36").Value = '  -1.75%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.31'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.76%  '
$ws.Range("E38").Value = '  -2.87%  '
$ws.Range("E39").Value = '  -4.00%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '50.65'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.19%  '
$ws.Range("E41").Value = '  -1.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '437.53'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.18%  '
$ws.Range("E43").Value = '  -1.64%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.23'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.47%  '
$ws.Range("E45").Value = '  +2.60%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0358'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.18%  '
$ws.Range("D47").Value = '2.791.01'
$ws.Range("E47").Value = '  -3.42%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '130.18'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.51%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.96'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.36%  '
$ws.Range("E51").Value = '  -0.12%  '
